# Add new PLC tag rows (14-23) to the "PLC Tags" sheet, following the
# exact same pattern used by the existing rows (2-13):
#   A = Name, B = "Marcas", C = "Bool", D = Logical Address,
#   E = "", F = "True", G = "True", H = "True", I = "", J = ""
#
# Row 2 is used as a template: copying it preserves the shared-string /
# text typing of the "Marcas", "Bool" and "True" values (plain text,
# not boolean) instead of Excel auto-converting the literal word "True"
# into a native boolean when it is typed directly into a cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PLC Tags")

$newRows = @(
    @("Salir",      "%M0.5"),
    @("A1_Manual",  "%M1000.0"),
    @("F1_Maunal",  "%M1000.1"),
    @("F2_E0",      "%M233.1"),
    @("F2_E1",      "%M233.2"),
    @("F2_E2",      "%M233.3"),
    @("F2_E3",      "%M233.4"),
    @("F2_E4",      "%M233.5"),
    @("F2_E5",      "%M233.6"),
    @("F2_E6",      "%M233.7")
)

$templateRow = $ws.Range("A2:J2")
$startRow = 14

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $name = $newRows[$i][0]
    $addr = $newRows[$i][1]

    $destRow = $ws.Range("A" + $r + ":J" + $r)
    $templateRow.Copy($destRow)

    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 4).Value = $addr
}
